# Update country data file: add MSME size-classification table (Oman)
# to the Summary sheet, pushing the existing MCI source-citation rows
# further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: the two "MCI" / "Ministry of Commerce..." source rows
# currently sitting at rows 25-26 need to end up at rows 31-32, so
# insert 6 blank rows right above them.
$ws.Range("A25:A30").EntireRow.Insert()

# --- New table header (row 22) ------------------------------------------------
$ws.Range("B22").Value = "Number of employees"
$ws.Range("C22").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D22").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B22:D22").Font.Bold = $true

# --- Micro row (23) ------------------------------------------------------------
$ws.Range("A23").Value = "Micro"
$ws.Range("B23").Value = "<5"
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = "< RO 25,000"

# --- Small row (24) -------------------------------------------------------------
$ws.Range("A24").Value = "Small"
$ws.Range("B24").Value = "5-9"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = "RO 25,000 < small < RO 250,000"

# --- Medium row (25) -------------------------------------------------------------
$ws.Range("A25").Value = "Medium"
$ws.Range("B25").Value = "10-99"
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = "RO 250,000 < small < RO 1,500,000"

# --- Large row (26) -------------------------------------------------------------
$ws.Range("A26").Value = "Large"
$ws.Range("B26").Value = ">99"
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = "> 1,500,000"

"done"
